# Project Sample Project is saved.TEST Author: admin. Type: SAVE.
#
# Row 11 of the Rules sheet relabels its "B11" cell from the shared string
# "R40" to the literal text "1".
#
# A plain   $ws.Range("B11").Value = "1"   lets Excel's normal type
# coercion turn a numeric-looking string into a real number, and forcing
# text via NumberFormat = "@" on B11 itself would mint a brand-new cell
# style for it (the real edit leaves B11's existing style untouched).
# So: stage the text "1" in a scratch cell formatted as Text, copy it, and
# paste only the VALUE into B11 (PasteSpecial xlPasteValues) -- this keeps
# B11's original style/formatting intact while still landing a genuine
# text value. The scratch cell is then fully cleared so it leaves no trace
# on the sheet (no stray style, no dimension growth).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$scratch = $ws.Range("Z1")
$scratch.NumberFormat = "@"
$scratch.Value = "1"

$scratch.Copy()
$ws.Range("B11").PasteSpecial(-4163)  # xlPasteValues

$scratch.Clear()
$excel.CutCopyMode = 0
